$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ballots")

$ws.Range("A31").Value = "Larry Stone"
$ws.Range("C31").Value = "x"
$ws.Range("D31").Value = "x"
$ws.Range("E31").Value = "x"
$ws.Range("F31").Value = "x"
$ws.Range("I31").Value = "x"
$ws.Range("J31").Value = "x"
$ws.Range("K31").Value = "x"
$ws.Range("O31").Value = "x"
$ws.Range("Q31").Value = "x"
$ws.Range("V31").Value = "x"
$ws.Range("AK31").Value = 10
$ws.Range("AL31").Value = "Twitter"
$ws.Range("AM31").Value = 43444
$ws.Range("AM30").Copy()
$ws.Range("AM31").PasteSpecial(-4122)

$ws.Range("A32").Value = "Paul White"
$ws.Range("E32").Value = "x"
$ws.Range("F32").Value = "x"
$ws.Range("I32").Value = "x"
$ws.Range("K32").Value = "x"
$ws.Range("O32").Value = "x"
$ws.Range("Q32").Value = "x"
$ws.Range("R32").Value = "x"
$ws.Range("U32").Value = "x"
$ws.Range("V32").Value = "x"
$ws.Range("AK32").Value = 9
$ws.Range("AL32").Value = "Email"
$ws.Range("AM32").Value = 43444
$ws.Range("AM30").Copy()
$ws.Range("AM32").PasteSpecial(-4122)

$ws.Range("AM28").Select()
